# Update "想去人数" (number of people interested) counts on the
# 展览 (Exhibitions) and 全部类型 (All Types) sheets.
#
# 展览 sheet: F2,F4,F6,F8,F9,F10,F11,F12,F13
# 全部类型 sheet: F2,F4,F6,F9,F10,F11,F12,F13,F14

$wb = $excel.ActiveWorkbook

$updates1 = @{
    2  = 80
    4  = 45
    6  = 112
    8  = 4478
    9  = 97
    10 = 4996
    11 = 564
    12 = 1261
    13 = 87
}

$updates4 = @{
    2  = 80
    4  = 45
    6  = 112
    9  = 4478
    10 = 97
    11 = 4996
    12 = 564
    13 = 1261
    14 = 87
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
